# Applies:
#  1) remove the "_GoBack" bookmark that currently sits right after the
#     "<trinh_ky>" placeholder.
#  2) split "Luật Quản lý thuế ngày 13 tháng 6 năm 2019 và các văn bản hướng
#     dẫn thi hành;" into four runs, replacing the date with the
#     "<luat_qlt_ngay>" placeholder, and insert a (collapsed) "_GoBack"
#     bookmark right before "và các văn bản...".
#
# (Part 1 has to run before part 2: "_GoBack" is a singleton bookmark name
#  as far as Bookmarks.Item/.Exists are concerned, so the old one has to be
#  gone before the new one is created, otherwise Item("_GoBack") keeps
#  resolving to whichever copy comes first in the document.)

$d = $word.ActiveDocument

# --- Part 1 : drop the "_GoBack" bookmark that used to follow <trinh_ky> ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

# --- Part 2 : rewrite the "Luật Quản lý thuế ..." sentence -----------------

$finder = $d.Content
$found = $finder.Find.Execute(
    "Luật Quản lý thuế ngày 13 tháng 6 năm 2019 và các văn bản hướng dẫn thi hành;",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # re-seat into a fresh Range: InsertXML misbehaves (appends instead of
    # replacing) on a Range that was itself the receiver of Find.Execute
    $target = $d.Range($finder.Start, $finder.End)

    $rPr = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr>'

    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r>' + $rPr + '<w:t xml:space="preserve">Luật Quản lý thuế </w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t>&lt;luat_qlt_ngay&gt;</w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r>' + $rPr + '<w:t>và các văn bản hướng dẫn thi hành;</w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $didXml = $false
    try {
        $target.InsertXML($xmlFrag)
        $didXml = $true
    } catch {
        $didXml = $false
    }

    if (-not $didXml) {
        # fallback: plain-text replace + re-add the bookmark at the seam
        # between "<luat_qlt_ngay> " and "và các văn bản..."
        $target2 = $d.Range($finder.Start, $finder.End)
        $target2.Text = "Luật Quản lý thuế <luat_qlt_ngay> và các văn bản hướng dẫn thi hành;"
        $seam = $finder.Start + ("Luật Quản lý thuế <luat_qlt_ngay> ").Length
        $bmRange = $d.Range($seam, $seam)
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}
